$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert 4 new rows above the current row 276 ("SessionData" block), shifting
# the existing rows 276-283 down to 280-287.
$ws.Rows("276:279").Insert()

# Expand the table (and its autofilter) to cover the 4 newly inserted rows.
$lo.Resize($ws.Range("A1:F287"))

# Populate the 4 new rows with the republished "sensitivity labels"
# snippets, writing cell-by-cell (instead of one bulk array assignment,
# which this host does not support for multi-cell ranges) and in the same
# order the original authoring session introduced each new shared string,
# so newly interned strings land at the same shared-string indices as the
# canonical workbook.
$ws.Cells.Item(279, 6).Value = "getSensitivityLabelsCatalogIsEnabled"
$ws.Cells.Item(278, 5).Value = "outlook-sensitivity-labels-sensitivity-labels-catalog"
$ws.Cells.Item(279, 3).Value = "getIsEnabledAsync"
$ws.Cells.Item(279, 2).Value = "SensitivityLabelsCatalog"
$ws.Cells.Item(278, 6).Value = "getSensitivityLabelsCatalog"
$ws.Cells.Item(277, 6).Value = "setSensitivityLabel"
$ws.Cells.Item(277, 2).Value = "SensitivityLabel"
$ws.Cells.Item(276, 6).Value = "getCurrentSensitivityLabel"
$ws.Cells.Item(276, 5).Value = "outlook-sensitivity-labels-sensitivity-label"

# Remaining cells: duplicate text already interned above, plus the
# "Package"/"Member ID" columns.
$ws.Cells.Item(279, 1).Value = "Office"
$ws.Cells.Item(279, 4).Value = 2
$ws.Cells.Item(279, 5).Value = "outlook-sensitivity-labels-sensitivity-labels-catalog"

$ws.Cells.Item(278, 1).Value = "Office"
$ws.Cells.Item(278, 2).Value = "SensitivityLabelsCatalog"
$ws.Cells.Item(278, 3).Value = "getAsync"
$ws.Cells.Item(278, 4).Value = 2

$ws.Cells.Item(277, 1).Value = "Office"
$ws.Cells.Item(277, 3).Value = "setAsync"
$ws.Cells.Item(277, 4).Value = 2
$ws.Cells.Item(277, 5).Value = "outlook-sensitivity-labels-sensitivity-label"

$ws.Cells.Item(276, 1).Value = "Office"
$ws.Cells.Item(276, 2).Value = "SensitivityLabel"
$ws.Cells.Item(276, 3).Value = "getAsync"
$ws.Cells.Item(276, 4).Value = 2

# Column E carries a column-wide right-aligned style; re-apply "Normal" on
# each new E cell so it relies on the column default instead of gaining an
# explicit (but redundant) style index, matching the rest of the table.
$ws.Cells.Item(276, 5).Style = "Normal"
$ws.Cells.Item(277, 5).Style = "Normal"
$ws.Cells.Item(278, 5).Style = "Normal"
$ws.Cells.Item(279, 5).Style = "Normal"

# These new rows render with a slightly different (custom) row height in the
# source workbook.
$ws.Rows("276:279").RowHeight = 14.25

# Restore the frozen-pane scroll anchor and move the active selection to
# match where the author was last working.
$ws.Range("F277").Select()
